# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the "time_taken" timestamps on the existing "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:20:37.908421"
$dataSheet.Range("F3").Value = "2021-10-05 14:20:37.908429"
$dataSheet.Range("F4").Value = "2021-10-05 14:20:37.908432"
$dataSheet.Range("F5").Value = "2021-10-05 14:20:37.908435"
$dataSheet.Range("F6").Value = "2021-10-05 14:20:37.908438"
$dataSheet.Range("F7").Value = "2021-10-05 14:20:37.908440"
$dataSheet.Range("F8").Value = "2021-10-05 14:20:37.908443"
$dataSheet.Range("F9").Value = "2021-10-05 14:20:37.908446"

# --- Add the new "metadata" sheet right after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (row 1, columns B-G)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row (row 2)
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Haematuria"
$metaSheet.Range("C2").Value = 99

# data_version needs to stay text "2.11" rather than being coerced to a number
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "2.11"
$metaSheet.Range("D2").Style = "Normal"

$metaSheet.Range("E2").Value = "2021-03-10T18:15:28.251141Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:37.904811"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/99/?format=json"

# Apply the bold/centered/bordered header style (same style used for headers on "data") to the header row
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)

$dataSheet.Range("B1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)

# Apply the same style used on "data"'s index column (A) to the metadata index cell
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$metaSheet.Range("A1").Select()
